# Refresh the ground-truth comparison results with the latest simulation run
# ("Run CircaDB analysis" / "Run CircadiPy simulations analysis").
# Rows 2-7 (columns A-G) are overwritten with the newly produced values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe keeps the textual labels (e.g. "01", "05", "10") as text
# instead of being reinterpreted as numbers, without touching number format
# or style (same General format + style index as before the edit).
$ws.Cells.Item(2, 1).Value = "'0"
$ws.Cells.Item(3, 1).Value = "'0.1"
$ws.Cells.Item(4, 1).Value = "'0.5"
$ws.Cells.Item(5, 1).Value = "'01"
$ws.Cells.Item(6, 1).Value = "'05"
$ws.Cells.Item(7, 1).Value = "'10"

# Row 2
$ws.Cells.Item(2, 2).Value = 0.06650000000031789
$ws.Cells.Item(2, 3).Value = 2.900790718740609 * [Math]::Pow(10, -13)
$ws.Cells.Item(2, 4).Value = 0.02250000000036732
$ws.Cells.Item(2, 5).Value = 0.2034520828107434
$ws.Cells.Item(2, 6).Value = 1.797942538808814 * [Math]::Pow(10, -13)
$ws.Cells.Item(2, 7).Value = 0.08717080933434017

# Row 3
$ws.Cells.Item(3, 2).Value = 0.2765000000003518
$ws.Cells.Item(3, 3).Value = 0.04899999999971758
$ws.Cells.Item(3, 4).Value = 0.1105000000003805
$ws.Cells.Item(3, 5).Value = 0.6360917779691767
$ws.Cells.Item(3, 6).Value = 0.2943450356299742
$ws.Cells.Item(3, 7).Value = 0.2713203825738432

# Row 4
$ws.Cells.Item(4, 2).Value = 0.04350000000031536
$ws.Cells.Item(4, 3).Value = 0.1160000000003066
$ws.Cells.Item(4, 4).Value = 0.1124999999996536
$ws.Cells.Item(4, 5).Value = 0.3921386872014265
$ws.Cells.Item(4, 6).Value = 0.2282191928826604
$ws.Cells.Item(4, 7).Value = 0.3466247971511013

# Row 5
$ws.Cells.Item(5, 2).Value = 0.03349999999969775
$ws.Cells.Item(5, 3).Value = 0.08300000000030323
$ws.Cells.Item(5, 4).Value = 0.03350000000036903
$ws.Cells.Item(5, 5).Value = 0.5612510579055742
$ws.Cells.Item(5, 6).Value = 0.4218897960368873
$ws.Cells.Item(5, 7).Value = 0.3819329129573403

# Row 6
$ws.Cells.Item(6, 2).Value = 0.1305000000003281
$ws.Cells.Item(6, 3).Value = 0.1184999999997293
$ws.Cells.Item(6, 4).Value = 0.01449999999963776
$ws.Cells.Item(6, 5).Value = 0.2242872042716541
$ws.Cells.Item(6, 6).Value = 0.2305054229297614
$ws.Cells.Item(6, 7).Value = 0.2802940420344944

# Row 7
$ws.Cells.Item(7, 2).Value = 0.1040000000003239
$ws.Cells.Item(7, 3).Value = 0.02699999999971432
$ws.Cells.Item(7, 4).Value = 0.02150000000036822
$ws.Cells.Item(7, 5).Value = 0.2290065501247217
$ws.Cells.Item(7, 6).Value = 0.05745432968892199
$ws.Cells.Item(7, 7).Value = 0.1453022711454006
